$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.705.43"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "3.149.89"
$ws.Range("E3").Value = "  +3.43%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'215.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").Value = "'626.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("E7").Value = "  +34.43%  "
$ws.Range("D8").Value = "'0.364"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.51%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "'0.804"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +21.65%  "
$ws.Range("B11").Value = "LidoStakedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D11").Value = "3.148.71"
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").Value = "'0.202"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.43%  "
$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").Value = "'5.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.45%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000243"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.86%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'35.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.89%  "
$ws.Range("D16").Value = "90.738.03"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "3.746.88"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "3.165.73"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").Value = "'3.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.65%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'14.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.63%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'465.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.32%  "
$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").Value = "'0.0000212"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("D23").Value = "'8.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.25%  "
$ws.Range("D24").Value = "'5.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.94%  "
$ws.Range("D25").Value = "'5.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.66%  "
$ws.Range("D26").Value = "'92.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.27%  "
$ws.Range("D27").Value = "'12.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.10%  "
$ws.Range("D28").Value = "3.334.80"
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'9.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +14.72%  "
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "'0.160"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.23%  "
$ws.Range("D33").Value = "'26.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +17.49%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.191"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +37.28%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'518.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.17%  "
$ws.Range("D36").Value = "'3.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("D37").Value = "'0.144"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.62%  "
$ws.Range("D38").Value = "'1.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.18%  "
$ws.Range("D39").Value = "'6.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.28%  "
$ws.Range("D40").Value = "'1.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.19%  "
$ws.Range("D41").Value = "'0.0897"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +30.98%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "'0.417"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.22%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.46%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'146.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "'4.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.36%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'44.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "'1.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.00%  "
$ws.Range("D51").Value = "'0.654"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.16%  "
